# Reorders the "Model Performance Table" rows so that, within each model
# group, the feature-set rows follow the canonical order
# (All Features, All Features PCA / Individual Features PCA / VGG / Top 3 Features)
# and the model groups are reordered as:
#   Logistic Regression, Random Forest, SVM, XGBoost
#
# This is implemented by rewriting each data row (rows 2-20) of the active
# sheet with the values that the commit moved into that row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = "Logistic Regression"
$ws.Cells.Item(2, 2).Value = "All Features"
$ws.Cells.Item(2, 3).Value = "C: 0.1, class_weight: balanced, l1_ratio: 1.0, multi_class: ovr, penalty: elasticnet, solver: saga"
$ws.Cells.Item(2, 4).Value = 0.9841745793769782
$ws.Cells.Item(2, 5).Value = 0.915708021093533
$ws.Cells.Item(2, 6).Value = 0.9106783075889859
$ws.Cells.Item(3, 1).Value = "Logistic Regression"
$ws.Cells.Item(3, 2).Value = "Individual Features PCA"
$ws.Cells.Item(3, 3).Value = "C: 0.01, class_weight: balanced, l1_ratio: 0.25, multi_class: ovr, penalty: elasticnet, solver: saga"
$ws.Cells.Item(3, 4).Value = 0.9580209895052474
$ws.Cells.Item(3, 5).Value = 0.9158759367194005
$ws.Cells.Item(3, 6).Value = 0.9093351242444594
$ws.Cells.Item(4, 1).Value = "Logistic Regression"
$ws.Cells.Item(4, 2).Value = "All Features PCA"
$ws.Cells.Item(4, 3).Value = "C: 0.1, class_weight: None, l1_ratio: 0.5, multi_class: ovr, penalty: elasticnet, solver: saga"
$ws.Cells.Item(4, 4).Value = 0.9921705813759787
$ws.Cells.Item(4, 5).Value = 0.9048806550097142
$ws.Cells.Item(4, 6).Value = 0.890530557421088
$ws.Cells.Item(5, 1).Value = "Logistic Regression"
$ws.Cells.Item(5, 2).Value = "VGG"
$ws.Cells.Item(5, 3).Value = "C: 0.1, class_weight: None, l1_ratio: 1.0, multi_class: multinomial, penalty: elasticnet, solver: saga"
$ws.Cells.Item(5, 4).Value = 0.9180409795102449
$ws.Cells.Item(5, 5).Value = 0.8635645295587011
$ws.Cells.Item(5, 6).Value = 0.857622565480188
$ws.Cells.Item(6, 1).Value = "Logistic Regression"
$ws.Cells.Item(6, 2).Value = "Top 3 Features"
$ws.Cells.Item(6, 3).Value = "C: 0.1, class_weight: None, l1_ratio: 1.0, multi_class: ovr, penalty: elasticnet, solver: saga"
$ws.Cells.Item(6, 4).Value = 0.9232050641345994
$ws.Cells.Item(6, 5).Value = 0.8580710519011934
$ws.Cells.Item(6, 6).Value = 0.8502350570852921
$ws.Cells.Item(7, 1).Value = "Random Forest"
$ws.Cells.Item(7, 2).Value = "All Features"
$ws.Cells.Item(7, 3).Value = "bootstrap: False, max_depth: 80, max_features: sqrt, min_samples_leaf: 1, n_estimators: 1500"
$ws.Cells.Item(7, 4).Value = 1
$ws.Cells.Item(7, 5).Value = 0.8615681376630585
$ws.Cells.Item(7, 6).Value = 0.8529214237743452
$ws.Cells.Item(8, 1).Value = "Random Forest"
$ws.Cells.Item(8, 2).Value = "Individual Features PCA"
$ws.Cells.Item(8, 3).Value = "bootstrap: False, max_depth: 40, max_features: sqrt, min_samples_leaf: 1, n_estimators: 500"
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 0.8084256175409381
$ws.Cells.Item(8, 6).Value = 0.7897918065815984
$ws.Cells.Item(9, 1).Value = "Random Forest"
$ws.Cells.Item(9, 2).Value = "VGG"
$ws.Cells.Item(9, 3).Value = "bootstrap: False, max_depth: 60, max_features: sqrt, min_samples_leaf: 2, n_estimators: 1500"
$ws.Cells.Item(9, 4).Value = 1
$ws.Cells.Item(9, 5).Value = 0.8062632528448516
$ws.Cells.Item(9, 6).Value = 0.7844190732034922
$ws.Cells.Item(10, 1).Value = "Random Forest"
$ws.Cells.Item(10, 2).Value = "Top 3 Features"
$ws.Cells.Item(10, 3).Value = "bootstrap: False, max_depth: 80, max_features: sqrt, min_samples_leaf: 2, n_estimators: 500"
$ws.Cells.Item(10, 4).Value = 1
$ws.Cells.Item(10, 5).Value = 0.7932697751873439
$ws.Cells.Item(10, 6).Value = 0.7575554063129617
$ws.Cells.Item(11, 1).Value = "Random Forest"
$ws.Cells.Item(11, 2).Value = "All Features PCA"
$ws.Cells.Item(11, 3).Value = "bootstrap: False, max_depth: 60, max_features: sqrt, min_samples_leaf: 1, n_estimators: 1500"
$ws.Cells.Item(11, 4).Value = 1
$ws.Cells.Item(11, 5).Value = 0.6989836247571468
$ws.Cells.Item(11, 6).Value = 0.7011417058428475
$ws.Cells.Item(12, 1).Value = "SVM"
$ws.Cells.Item(12, 2).Value = "Individual Features PCA"
$ws.Cells.Item(12, 3).Value = "C: 10.0, class_weight: balanced, gamma: 1e-05, kernel: rbf"
$ws.Cells.Item(12, 4).Value = 0.9631850741296019
$ws.Cells.Item(12, 5).Value = 0.904880932556203
$ws.Cells.Item(12, 6).Value = 0.9046339825386165
$ws.Cells.Item(13, 1).Value = "SVM"
$ws.Cells.Item(13, 2).Value = "All Features PCA"
$ws.Cells.Item(13, 3).Value = "C: 10.0, class_weight: balanced, gamma: 1e-05, kernel: rbf"
$ws.Cells.Item(13, 4).Value = 0.9643511577544561
$ws.Cells.Item(13, 5).Value = 0.8827245351096309
$ws.Cells.Item(13, 6).Value = 0.8985896574882472
$ws.Cells.Item(14, 1).Value = "SVM"
$ws.Cells.Item(14, 2).Value = "VGG"
$ws.Cells.Item(14, 3).Value = "C: 10.0, class_weight: balanced, gamma: scale, kernel: rbf"
$ws.Cells.Item(14, 4).Value = 1
$ws.Cells.Item(14, 5).Value = 0.8697293921731891
$ws.Cells.Item(14, 6).Value = 0.8623237071860309
$ws.Cells.Item(15, 1).Value = "SVM"
$ws.Cells.Item(15, 2).Value = "Top 3 Features"
$ws.Cells.Item(15, 3).Value = "C: 10.0, class_weight: balanced, gamma: 0.00025118864315095795, kernel: rbf"
$ws.Cells.Item(15, 4).Value = 0.9890054972513743
$ws.Cells.Item(15, 5).Value = 0.8489078545656399
$ws.Cells.Item(15, 6).Value = 0.8388179986568166
$ws.Cells.Item(16, 1).Value = "XGBoost"
$ws.Cells.Item(16, 2).Value = "All Features"
$ws.Cells.Item(16, 3).Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 1, n_estimators: 300"
$ws.Cells.Item(16, 4).Value = 1
$ws.Cells.Item(16, 5).Value = 0.9085460727171804
$ws.Cells.Item(16, 6).Value = 0.9012760241773002
$ws.Cells.Item(17, 1).Value = "XGBoost"
$ws.Cells.Item(17, 2).Value = "Individual Features PCA"
$ws.Cells.Item(17, 3).Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 3, n_estimators: 300"
$ws.Cells.Item(17, 4).Value = 1
$ws.Cells.Item(17, 5).Value = 0.890555925617541
$ws.Cells.Item(17, 6).Value = 0.8757555406312961
$ws.Cells.Item(18, 1).Value = "XGBoost"
$ws.Cells.Item(18, 2).Value = "All Features PCA"
$ws.Cells.Item(18, 3).Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 3, n_estimators: 300"
$ws.Cells.Item(18, 4).Value = 1
$ws.Cells.Item(18, 5).Value = 0.8604021648626144
$ws.Cells.Item(18, 6).Value = 0.8643384822028207
$ws.Cells.Item(19, 1).Value = "XGBoost"
$ws.Cells.Item(19, 2).Value = "VGG"
$ws.Cells.Item(19, 3).Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 1, n_estimators: 300"
$ws.Cells.Item(19, 4).Value = 1
$ws.Cells.Item(19, 5).Value = 0.8637334165972801
$ws.Cells.Item(19, 6).Value = 0.8542646071188718
$ws.Cells.Item(20, 1).Value = "XGBoost"
$ws.Cells.Item(20, 2).Value = "Top 3 Features"
$ws.Cells.Item(20, 3).Value = "learning_rate: 0.5, max_depth: 3, min_child_weight: 1, n_estimators: 300"
$ws.Cells.Item(20, 4).Value = 1
$ws.Cells.Item(20, 5).Value = 0.8597356369691923
$ws.Cells.Item(20, 6).Value = 0.8488918737407656
